$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 679
$ws.Range("F6").Value = 603
$ws.Range("F10").Value = 6077
$ws.Range("F11").Value = 658
$ws.Range("F12").Value = 1063
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 259
$ws.Range("F17").Value = 594
$ws.Range("F18").Value = 1020
$ws.Range("F19").Value = 52
$ws.Range("F21").Value = 193
$ws.Range("F22").Value = 1352
$ws.Range("F24").Value = 1026
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 2090
$ws.Range("F27").Value = 200
$ws.Range("F28").Value = 36
$ws.Range("F29").Value = 371
$ws.Range("F31").Value = 3375

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 29
$ws.Range("F5").Value = 348
$ws.Range("F6").Value = 103
$ws.Range("F9").Value = 20
$ws.Range("F10").Value = 658
$ws.Range("F15").Value = 83
$ws.Range("F16").Value = 631
$ws.Range("F18").Value = 68
$ws.Range("F20").Value = 363
$ws.Range("F21").Value = 302
$ws.Range("F22").Value = 4064
$ws.Range("F28").Value = 76
$ws.Range("F34").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1158
$ws.Range("F8").Value = 1514
$ws.Range("F10").Value = 120
$ws.Range("F12").Value = 683

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 29
$ws.Range("F6").Value = 1158
$ws.Range("F7").Value = 1514
$ws.Range("F9").Value = 120
$ws.Range("F10").Value = 679
$ws.Range("F12").Value = 683
$ws.Range("F13").Value = 603
$ws.Range("F14").Value = 103
$ws.Range("F18").Value = 6077
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 658
$ws.Range("F21").Value = 1063
$ws.Range("F22").Value = 658
$ws.Range("F23").Value = 259
$ws.Range("F26").Value = 594
$ws.Range("F28").Value = 83
$ws.Range("F30").Value = 68
$ws.Range("F31").Value = 1020
$ws.Range("F34").Value = 363
$ws.Range("F40").Value = 73
$ws.Range("F41").Value = 76
$ws.Range("F43").Value = 2090
$ws.Range("F46").Value = 200
$ws.Range("F47").Value = 36
$ws.Range("F48").Value = 371
$ws.Range("F49").Value = 3375
